# Fruta / hortaliza, semanal
# A new weekly record (Fecha=44705 / 2022-05-24, Volumen=120) is inserted
# into the "Ciboulette" price table, as row 157. All subsequent rows
# (old 157..316) shift down by one (to 158..317), extending the used
# range from A1:R316 to A1:R317.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157; existing rows 157-316 shift down to 158-317.
$ws.Rows(157).Insert()

# Populate the newly inserted row 157 with the new record.
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = "Femacal de La Calera"
$ws.Range("C157").Value = "Coquimbo"
$ws.Range("D157").Value = 44705
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = 100112039
$ws.Range("G157").Value = "Ciboulette"
$ws.Range("H157").Value = "Sin especificar"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 120
$ws.Range("K157").Value = 1500
$ws.Range("L157").Value = 1500
$ws.Range("M157").Value = 1500
$ws.Range("N157").Value = "$/docena de atados"
$ws.Range("O157").Value = "Provincia de Quillota"
$ws.Range("P157").Value = 500
$ws.Range("Q157").Value = 3
$ws.Range("R157").Value = "Hortaliza"
